$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string (rich text run) edits: volume number and report week dates ---
$ws.Range("A8").Characters(21, 2).Text = "45"
$ws.Range("C9").Characters(48, 9).Text = "11/10/2024"
$ws.Range("C9").Characters(27, 10).Text = "11/4/2024"

# --- Plain value updates (style unchanged) ---
$ws.Range("N14").Value = -57.142857142857
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 10.526315789473
$ws.Range("N15").Value = -22.222222222222
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = -4.545454545454
$ws.Range("L16").Value = -19.230769230769
$ws.Range("M16").Value = -39.506172839506
$ws.Range("N16").Value = -85.588235294117
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -43.75
$ws.Range("I17").Value = 277
$ws.Range("J17").Value = 339
$ws.Range("K17").Value = -18.289085545722
$ws.Range("L17").Value = -15.805471124620
$ws.Range("M17").Value = 38.5
$ws.Range("N17").Value = -13.975155279503
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = -28.813559322033
$ws.Range("L18").Value = -36.842105263157
$ws.Range("M18").Value = -72.185430463576
$ws.Range("N18").Value = -93.285371702637
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 26
$ws.Range("I19").Value = 307
$ws.Range("J19").Value = 292
$ws.Range("K19").Value = 5.136986301369
$ws.Range("L19").Value = -6.116207951070
$ws.Range("M19").Value = 1.320132013201
$ws.Range("N19").Value = -38.476953907815
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = -70.588235294117
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = -45.945945945945
$ws.Range("I20").Value = 254
$ws.Range("J20").Value = 271
$ws.Range("K20").Value = -6.273062730627
$ws.Range("L20").Value = 13.901345291479
$ws.Range("M20").Value = -0.392156862745
$ws.Range("N20").Value = -91.522029372496
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -56.097560975609
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -32.456140350877
$ws.Range("I21").Value = 1093
$ws.Range("J21").Value = 1195
$ws.Range("K21").Value = -8.535564853556
$ws.Range("L21").Value = -10.262725779967
$ws.Range("M21").Value = -17.134192570128
$ws.Range("N21").Value = -82.146357399542
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 77.777777777777
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = 51.315789473684
$ws.Range("I24").Value = 1227
$ws.Range("J24").Value = 1200
$ws.Range("K24").Value = 2.25
$ws.Range("L24").Value = -13.773717498243
$ws.Range("M24").Value = 76.293103448275
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 144.444444444444
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 181.481481481481
$ws.Range("I25").Value = 634
$ws.Range("J25").Value = 405
$ws.Range("K25").Value = 56.543209876543
$ws.Range("L25").Value = 11.228070175438
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 9.090909090909
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 9.090909090909
$ws.Range("I26").Value = 552
$ws.Range("J26").Value = 483
$ws.Range("K26").Value = 14.285714285714
$ws.Range("L26").Value = 15.240083507306
$ws.Range("M26").Value = 0.181488203266
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 23.333333333333
$ws.Range("L27").Value = 15.625
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 47
$ws.Range("K28").Value = -2.083333333333
$ws.Range("L28").Value = 4.444444444444
$ws.Range("N29").Value = -53.571428571428
$ws.Range("N30").Value = -60.869565217391

# --- Cells whose content type changes (number <-> text placeholder): set value, then fix style ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("C16").Value = "'0"
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("C28").Value = 3

# Reference cells that keep style 13 (text/General), 14 (#,##0) and 15 (#,##0.0) throughout
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
